$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell's displayed text on the crypto listing sheet
# (Price in column D, 1h Volume % change in column E). Cells flagged
# Text = $true hold numeric-looking strings (e.g. "485.81") that must stay
# plain text (matching the sheet's existing inline-string cells) instead of
# being auto-converted to numbers by Excel's input parser, so we briefly
# switch the cell to the Text number format while assigning, then restore
# the Normal style so formatting is left exactly as it was found.
$updates = @(
    @{ Cell = "D2"; Value = "53.652.47"; Text = $false },
    @{ Cell = "E2"; Value = "  -5.20%  "; Text = $false },
    @{ Cell = "D3"; Value = "2.206.25"; Text = $false },
    @{ Cell = "E3"; Value = "  -7.60%  "; Text = $false },
    @{ Cell = "E4"; Value = "  -0.08%  "; Text = $false },
    @{ Cell = "D5"; Value = "485.81"; Text = $true },
    @{ Cell = "E5"; Value = "  -4.29%  "; Text = $false },
    @{ Cell = "D6"; Value = "125.20"; Text = $true },
    @{ Cell = "E6"; Value = "  -4.45%  "; Text = $false },
    @{ Cell = "D7"; Value = "0.995"; Text = $true },
    @{ Cell = "E7"; Value = "  -0.21%  "; Text = $false },
    @{ Cell = "D8"; Value = "0.520"; Text = $true },
    @{ Cell = "E8"; Value = "  -4.83%  "; Text = $false },
    @{ Cell = "D9"; Value = "2.232.46"; Text = $false },
    @{ Cell = "E9"; Value = "  -7.04%  "; Text = $false },
    @{ Cell = "D10"; Value = "0.0923"; Text = $true },
    @{ Cell = "E10"; Value = "  -6.96%  "; Text = $false },
    @{ Cell = "E11"; Value = "  -0.74%  "; Text = $false },
    @{ Cell = "D12"; Value = "0.317"; Text = $true },
    @{ Cell = "E12"; Value = "  -3.34%  "; Text = $false },
    @{ Cell = "D13"; Value = "4.62"; Text = $true },
    @{ Cell = "E13"; Value = "  -5.05%  "; Text = $false },
    @{ Cell = "D14"; Value = "2.600.38"; Text = $false },
    @{ Cell = "E14"; Value = "  -7.58%  "; Text = $false },
    @{ Cell = "D15"; Value = "21.19"; Text = $true },
    @{ Cell = "E15"; Value = "  -2.33%  "; Text = $false },
    @{ Cell = "D16"; Value = "53.583.44"; Text = $false },
    @{ Cell = "E16"; Value = "  -5.28%  "; Text = $false },
    @{ Cell = "D17"; Value = "0.0000128"; Text = $true },
    @{ Cell = "E17"; Value = "  -4.67%  "; Text = $false },
    @{ Cell = "D18"; Value = "2.227.90"; Text = $false },
    @{ Cell = "E18"; Value = "  -7.61%  "; Text = $false },
    @{ Cell = "D19"; Value = "9.64"; Text = $true },
    @{ Cell = "E19"; Value = "  -5.12%  "; Text = $false },
    @{ Cell = "D20"; Value = "3.97"; Text = $true },
    @{ Cell = "E20"; Value = "  -2.13%  "; Text = $false },
    @{ Cell = "D21"; Value = "295.36"; Text = $true },
    @{ Cell = "E21"; Value = "  -4.96%  "; Text = $false },
    @{ Cell = "D22"; Value = "6.15"; Text = $true },
    @{ Cell = "E22"; Value = "  -2.70%  "; Text = $false },
    @{ Cell = "D23"; Value = "0.997"; Text = $true },
    @{ Cell = "E23"; Value = "  -0.24%  "; Text = $false },
    @{ Cell = "D24"; Value = "63.20"; Text = $true },
    @{ Cell = "E24"; Value = "  -4.57%  "; Text = $false },
    @{ Cell = "D25"; Value = "0.996"; Text = $true },
    @{ Cell = "E25"; Value = "  +0.15%  "; Text = $false },
    @{ Cell = "D26"; Value = "0.366"; Text = $true },
    @{ Cell = "E26"; Value = "  -1.30%  "; Text = $false },
    @{ Cell = "D27"; Value = "0.147"; Text = $true },
    @{ Cell = "E27"; Value = "  -0.32%  "; Text = $false },
    @{ Cell = "D28"; Value = "2.308.16"; Text = $false },
    @{ Cell = "E28"; Value = "  -7.59%  "; Text = $false },
    @{ Cell = "D29"; Value = "7.05"; Text = $true },
    @{ Cell = "E29"; Value = "  -3.58%  "; Text = $false },
    @{ Cell = "D30"; Value = "163.32"; Text = $true },
    @{ Cell = "E30"; Value = "  -5.61%  "; Text = $false },
    @{ Cell = "D31"; Value = "1.58"; Text = $true },
    @{ Cell = "E31"; Value = "  -4.54%  "; Text = $false },
    @{ Cell = "D32"; Value = "0.998"; Text = $true },
    @{ Cell = "E32"; Value = "  -0.14%  "; Text = $false },
    @{ Cell = "D33"; Value = "0.0₃0669"; Text = $false },
    @{ Cell = "E33"; Value = "  -6.93%  "; Text = $false },
    @{ Cell = "D34"; Value = "5.77"; Text = $true },
    @{ Cell = "E34"; Value = "  -1.87%  "; Text = $false },
    @{ Cell = "D35"; Value = "0.993"; Text = $true },
    @{ Cell = "E35"; Value = "  -0.25%  "; Text = $false },
    @{ Cell = "D36"; Value = "1.05"; Text = $true },
    @{ Cell = "E36"; Value = "  -3.88%  "; Text = $false },
    @{ Cell = "D37"; Value = "17.33"; Text = $true },
    @{ Cell = "E37"; Value = "  -2.45%  "; Text = $false },
    @{ Cell = "D38"; Value = "1.17"; Text = $true },
    @{ Cell = "E38"; Value = "  -1.98%  "; Text = $false },
    @{ Cell = "D39"; Value = "0.843"; Text = $true },
    @{ Cell = "E39"; Value = "  +1.43%  "; Text = $false },
    @{ Cell = "E40"; Value = "  -5.53%  "; Text = $false },
    @{ Cell = "D41"; Value = "35.18"; Text = $true },
    @{ Cell = "E41"; Value = "  -3.77%  "; Text = $false },
    @{ Cell = "D42"; Value = "0.367"; Text = $true },
    @{ Cell = "E42"; Value = "  -1.04%  "; Text = $false },
    @{ Cell = "D43"; Value = "1.37"; Text = $true },
    @{ Cell = "E43"; Value = "  -2.15%  "; Text = $false },
    @{ Cell = "D44"; Value = "3.29"; Text = $true },
    @{ Cell = "E44"; Value = "  -3.39%  "; Text = $false },
    @{ Cell = "D45"; Value = "126.23"; Text = $true },
    @{ Cell = "E45"; Value = "  -2.55%  "; Text = $false },
    @{ Cell = "D46"; Value = "4.84"; Text = $true },
    @{ Cell = "E46"; Value = "  +1.81%  "; Text = $false },
    @{ Cell = "D47"; Value = "0.0881"; Text = $true },
    @{ Cell = "E47"; Value = "  -2.09%  "; Text = $false },
    @{ Cell = "D48"; Value = "0.534"; Text = $true },
    @{ Cell = "E48"; Value = "  -6.82%  "; Text = $false },
    @{ Cell = "D49"; Value = "233.88"; Text = $true },
    @{ Cell = "E49"; Value = "  -3.56%  "; Text = $false },
    @{ Cell = "D50"; Value = "0.0472"; Text = $true },
    @{ Cell = "E50"; Value = "  -2.57%  "; Text = $false },
    @{ Cell = "E51"; Value = "  -3.88%  "; Text = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.Text) {
        $r.NumberFormat = "@"
        $r.Value = $u.Value
        $r.Style = "Normal"
    } else {
        $r.Value = $u.Value
    }
}
